$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift names in column A up by one (rows 3-14)
$ws.Range("A3").Value  = "bellingham"
$ws.Range("A4").Value  = "benzema"
$ws.Range("A5").Value  = "carvajal"
$ws.Range("A6").Value  = "courtois"
$ws.Range("A7").Value  = "kevin"
$ws.Range("A8").Value  = "mbappe"
$ws.Range("A9").Value  = "messi"
$ws.Range("A10").Value = "modric"
$ws.Range("A11").Value = "mtp"
$ws.Range("A12").Value = "neymar"
$ws.Range("A13").Value = "ronaldo"
$ws.Range("A14").Value = "vietANH"

# Append new attendance rows 16 and 17
$ws.Range("A16").Value = "QA"
$ws.Range("B16").Value = "Có mặt"
$ws.Range("C16").Value = "16:09:43"
$ws.Range("D16").Value = "N/A"
$ws.Range("E16").Value = "Chưa checkout"

$ws.Range("A17").Value = "anhtonton"
$ws.Range("B17").Value = "Có mặt"
$ws.Range("C17").Value = "16:09:37"
$ws.Range("D17").Value = "N/A"
$ws.Range("E17").Value = "Chưa checkout"
